$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "URL of site evaluated:" -> "URL of site evaluated " (drop colon, keep
#    trailing space) followed by the actual URL and a shorter underscore run,
#    replacing the old "_________________________" placeholder.
# ---------------------------------------------------------------------------
$rngUrlLabel = $d.Content
$rngUrlLabel.Find.Text = "URL of site evaluated:"
$rngUrlLabel.Find.Execute() | Out-Null
$rngUrlLabel.Text = "URL of site evaluated "
# Touch the font explicitly so this stays its own run (and picks up
# xml:space="preserve" for the trailing space) instead of silently merging
# back into whatever follows it.
$rngUrlLabel.Font.Name = "Arial"
$rngUrlLabel.Font.NameBi = "Arial"

# Replace the placeholder underscores with the URL immediately followed by a
# shorter run of underscores, all written as one go so the run inherits the
# exact formatting (Arial / szCs 20) of the text it replaces.
$rngUrl = $d.Content
$rngUrl.Find.Text = "_________________________"
$rngUrl.Find.Execute() | Out-Null
$rngUrl.Text = "https://jgunasingham.github.io/islt_7310/index.html______________________"
$rngUrl.Font.Name = "Arial"
$rngUrl.Font.NameBi = "Arial"

# Split the trailing underscores back out into their own run.
$rngUrlTail = $d.Content
$rngUrlTail.Find.Text = "______________________"
$rngUrlTail.Find.Execute() | Out-Null
$rngUrlTail.Font.Name = "Arial"
$rngUrlTail.Font.NameBi = "Arial"

# The empty paragraph that used to sit directly below the URL line is
# removed (its paragraph mark is deleted, merging it away).
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -eq [char]13 -or $p.Range.Text -eq "") {
        $prev = $paras.Item($i - 1)
        if ($prev.Range.Text -like "URL of site evaluated*") {
            $p.Range.Delete()
            break
        }
    }
}

# ---------------------------------------------------------------------------
# 2) "Author of site evaluated:  _________________________________" ->
#    "Author of site evaluated:  Jonathan Gunasingham" (name split across its
#    own runs, "Gunasingham" wrapped in a spelling proof-error pair in real
#    Word - not reproducible through this object model).
# ---------------------------------------------------------------------------
$rngAuthorBlank = $d.Content
$rngAuthorBlank.Find.Text = "  _________________________________"
$rngAuthorBlank.Find.Execute() | Out-Null
$rngAuthorBlank.Text = "  Jonathan Gunasingham"
$rngAuthorBlank.Font.Name = "Arial"
$rngAuthorBlank.Font.NameBi = "Arial"

$rngJonathan = $d.Content
$rngJonathan.Find.Text = "Jonathan "
$rngJonathan.Find.Execute() | Out-Null
$rngJonathan.Font.Name = "Arial"
$rngJonathan.Font.NameBi = "Arial"

$rngGuna = $d.Content
$rngGuna.Find.Text = "Gunasingham"
$rngGuna.Find.Execute() | Out-Null
$rngGuna.Font.Name = "Arial"
$rngGuna.Font.NameBi = "Arial"

# ---------------------------------------------------------------------------
# 3) Drop a "_GoBack" bookmark (the marker Word leaves at the last edit
#    point when a document is saved) into the blank paragraph two
#    paragraphs below the author line.
# ---------------------------------------------------------------------------
$paras = $d.Paragraphs
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -eq [char]13) {
        $prev = $paras.Item($i - 1)
        if ($prev.Range.Text -eq [char]13) {
            $prevPrev = $paras.Item($i - 2)
            if ($prevPrev.Range.Text -like "Author of site evaluated*") {
                $bmRange = $p.Range
                $bmRange.Collapse(1)
                $d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
                break
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 4) Footer page-number field: cached result "3" -> "1".
# ---------------------------------------------------------------------------
$footer = $d.Sections.Item(1).Footers.Item(1)
$rngPage = $footer.Range
$rngPage.Find.Execute("3", $true, $false, $false, $false, $false, $true, 1, $false, "1", 2) | Out-Null

Write-Output "edit complete"
